# Added function to allow changing df in excel
# Reworks the setup sheet's header/data row: adds an "index" column after
# "datnumplus", reorders "datetime" after it, and drops the trailing
# fd0adc/fd1adc/fd2adc/fd3adc columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) -------------------------------------------------
$ws.Range("A1").Value = "datnumplus"
$ws.Range("B1").Value = "index"
$ws.Range("C1").Value = "datetime"
$ws.Range("D1").Value = "i_sense"
$ws.Range("E1").Value = "FastScan"
$ws.Range("F1").Value = "FastScanCh0"
$ws.Range("G1").Value = "FastScanCh1"
$ws.Range("H1").Value = "FastScanCh2"
$ws.Range("I1").Value = "FastScanCh3"

# --- Row 2 (sample data) ----------------------------------------------
# A2 keeps its existing (header) style but becomes a plain number.
$ws.Range("A2").Value = 0

# B2 is a brand new data cell with no special formatting.
$ws.Range("B2").ClearFormats() | Out-Null
$ws.Range("B2").Value = 0

# C2 takes over the datetime string that used to live in A2.
$ws.Range("C2").Value = "Wednesday, January 1, 2020 00:00:00"

# D2:I2 were already blank placeholders and stay that way untouched.

# --- Drop the old fd0adc..fd3adc columns (now J:L after the shuffle) --
$ws.Range("J1:L2").Delete(-4159) | Out-Null
